$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) C40: "Expedicion" -> "Produccion - Expedicion" (new label for the
#    Production/Expedition building). Style (s=1) already correct.
# ---------------------------------------------------------------------------
$ws.Range("C40").Value = "Produccion - Expedici" + [char]0x00F3 + "n"

# ---------------------------------------------------------------------------
# 2) Row 44 was "Acces point" -> becomes "Switch" (style s=4 unchanged).
# ---------------------------------------------------------------------------
$ws.Range("C44").Value = "Switch"

# ---------------------------------------------------------------------------
# 3) Row 45 used to be a blank filler row; it now carries the
#    "Acces point" line (qty 1 x 35W) that used to live in row 44.
#    Styles (s=4/7/7/7) already match a blank row, so no format paste
#    is required there.
# ---------------------------------------------------------------------------
$ws.Range("C45").Value = "Acces point"
$ws.Range("D45").Value = 1
$ws.Range("E45").Value = 35
$ws.Range("F45").Formula = "=D45*E45"

# ---------------------------------------------------------------------------
# 4) Row 46 becomes the new blank filler row (C/D/E/F, style s=4/7/7/7).
#    Grab the format from the existing blank-row pattern (row 45 before
#    we touched it is gone now, so copy from the still-blank D/E/F donor
#    cells elsewhere, and from C44 for the C column).
# ---------------------------------------------------------------------------
$ws.Range("C44").Copy()
$ws.Range("C46").PasteSpecial(-4122)
$ws.Range("D44").Copy()
$ws.Range("D46").PasteSpecial(-4122)
$ws.Range("E44").Copy()
$ws.Range("E46").PasteSpecial(-4122)
$ws.Range("F44").Copy()
$ws.Range("F46").PasteSpecial(-4122)

$ws.Range("E46").ClearContents()
$ws.Range("F46").ClearContents()
$ws.Range("K46").Clear()

# ---------------------------------------------------------------------------
# 5) Row 47 becomes the "Total:" row (was row 46's job).
# ---------------------------------------------------------------------------
$ws.Range("E13").Copy()
$ws.Range("E47").PasteSpecial(-4122)
$ws.Range("F44").Copy()
$ws.Range("F47").PasteSpecial(-4122)

$ws.Range("E47").Value = "Total:"
$ws.Range("F47").Formula = "=SUM(F42:F46)"

$ws.Range("H47").Clear()
$ws.Range("I47").Clear()
$ws.Range("K47").Clear()
$ws.Range("M47").Clear()

# ---------------------------------------------------------------------------
# 6) Row 48 becomes the "C/Resguardo" row (was row 47's job).
# ---------------------------------------------------------------------------
$ws.Range("E13").Copy()
$ws.Range("E48").PasteSpecial(-4122)
$ws.Range("F13").Copy()
$ws.Range("F48").PasteSpecial(-4122)

$ws.Range("E48").Value = "C/Resguardo"
$ws.Range("F48").Formula = "=F47+(F47*$I$41)"

$ws.Range("H48").Clear()
$ws.Range("I48").Clear()

# ---------------------------------------------------------------------------
# 7) Row 49 becomes the "A" row (was row 48's job). It did not exist
#    before, so paste formats from donor cells first.
# ---------------------------------------------------------------------------
$ws.Range("C44").Copy()
$ws.Range("E49").PasteSpecial(-4122)
$ws.Range("F14").Copy()
$ws.Range("F49").PasteSpecial(-4122)

$ws.Range("E49").Value = "A"
$ws.Range("F49").Formula = "=F48/$L$41"

# ---------------------------------------------------------------------------
# 8) I42's consumption-total formula referenced F47 (the old
#    "C/Resguardo" subtotal); it must now reference F48, the cell that
#    now plays that role. F57/L57 stay untouched.
# ---------------------------------------------------------------------------
$ws.Range("I42").Formula = "=F13+L10+F23+L23+F35+L35+F48+F57+L57"

$excel.Calculate()
